# Auto-generated edit script applying cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.193.19"
$ws.Range("E2").Value = "  -4.51%  "
$ws.Range("D3").Value = "1.654.88"
$ws.Range("E3").Value = "  -3.85%  "
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'216.05"
$ws.Range("E5").Value = "  -4.37%  "
$ws.Range("D6").Value = "'0.5134"
$ws.Range("E6").Value = "  -3.49%  "
$ws.Range("D7").Value = "'1.005"
$ws.Range("E7").Value = "  +0.30%  "
$ws.Range("D8").Value = "'0.2600"
$ws.Range("E8").Value = "  -2.67%  "
$ws.Range("D9").Value = "'0.06449"
$ws.Range("E9").Value = "  -4.09%  "
$ws.Range("D10").Value = "'19.97"
$ws.Range("E10").Value = "  -5.19%  "
$ws.Range("D11").Value = "'0.07783"
$ws.Range("E11").Value = "  +0.79%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.655.15"
$ws.Range("E12").Value = "  -3.69%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'4.299"
$ws.Range("E13").Value = "  -4.92%  "
$ws.Range("D14").Value = "1.881.50"
$ws.Range("E14").Value = "  -3.87%  "
$ws.Range("D15").Value = "'0.5535"
$ws.Range("E16").Value = "  -2.67%  "
$ws.Range("D17").Value = "'64.25"
$ws.Range("E17").Value = "  -5.95%  "
$ws.Range("D18").Value = "26.178.73"
$ws.Range("E18").Value = "  -4.64%  "
$ws.Range("D19").Value = "'1.006"
$ws.Range("E19").Value = "  +0.27%  "
$ws.Range("D20").Value = "'211.04"
$ws.Range("E20").Value = "  -6.37%  "
$ws.Range("D21").Value = "'4.412"
$ws.Range("E21").Value = "  -5.76%  "
$ws.Range("E22").Value = "  -4.34%  "
$ws.Range("D23").Value = "'6.041"
$ws.Range("E23").Value = "  -0.43%  "
$ws.Range("E24").Value = "  +0.20%  "
$ws.Range("D25").Value = "'143.96"
$ws.Range("E25").Value = "  -0.28%  "
$ws.Range("E26").Value = "  +2.67%  "
$ws.Range("D27").Value = "'0.1176"
$ws.Range("E27").Value = "  -3.10%  "
$ws.Range("D28").Value = "'6.988"
$ws.Range("E28").Value = "  -3.98%  "
$ws.Range("D29").Value = "'15.85"
$ws.Range("E29").Value = "  -3.34%  "
$ws.Range("D30").Value = "'0.05105"
$ws.Range("E30").Value = "  -5.35%  "
$ws.Range("D31").Value = "'1.243"
$ws.Range("E31").Value = "  -4.18%  "
$ws.Range("D32").Value = "'3.382"
$ws.Range("E32").Value = "  -3.24%  "
$ws.Range("D33").Value = "'3.230"
$ws.Range("E33").Value = "  -5.99%  "
$ws.Range("D34").Value = "'1.566"
$ws.Range("E34").Value = "  -4.15%  "
$ws.Range("D35").Value = "'2.736"
$ws.Range("E35").Value = "  -4.68%  "
$ws.Range("D36").Value = "'0.9256"
$ws.Range("E36").Value = "  -3.49%  "
$ws.Range("D37").Value = "'2.349"
$ws.Range("E37").Value = "  -1.66%  "
$ws.Range("D38").Value = "'0.5736"
$ws.Range("E38").Value = "  -2.98%  "
$ws.Range("D39").Value = "1.166.70"
$ws.Range("E39").Value = "  +1.05%  "
$ws.Range("D40").Value = "'0.01588"
$ws.Range("E40").Value = "  -4.32%  "
$ws.Range("D41").Value = "'2.557"
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("E42").Value = "  +0.25%  "
$ws.Range("D43").Value = "'5.668"
$ws.Range("E43").Value = "  -2.35%  "
$ws.Range("D44").Value = "'0.8232"
$ws.Range("E44").Value = "  -2.58%  "
$ws.Range("D45").Value = "'100.28"
$ws.Range("E45").Value = "  -0.70%  "
$ws.Range("D46").Value = "1.791.21"
$ws.Range("E46").Value = "  -3.86%  "
$ws.Range("E47").Value = "  +3.75%  "
$ws.Range("D48").Value = "'0.4547"
$ws.Range("E48").Value = "  -0.78%  "
$ws.Range("D49").Value = "'55.60"
$ws.Range("E49").Value = "  -4.22%  "
$ws.Range("D50").Value = "'1.008"
$ws.Range("E50").Value = "  +0.15%  "
$ws.Range("D51").Value = "'7.867"
$ws.Range("E51").Value = "  -3.49%  "
